# Apply the changes described by the diff:
# 1. Update the email value (and shared string) on InfoBasica!A2
# 2. Make InfoBasica the active sheet with A4 selected (removes tabSelected
#    from BuyBook and activeTab from the workbook view)

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("InfoBasica")
$wsBuy  = $wb.Worksheets.Item("BuyBook")

# Update the e-mail address. Both InfoBasica!A2 and BuyBook!A2 point at the
# same shared string, so update both cells to keep them sharing one entry.
$wsInfo.Range("A2").Value = "pruebareto1111@yopmail.com"
$wsBuy.Range("A2").Value  = "pruebareto1111@yopmail.com"

# Activate InfoBasica and select cell A4, leaving BuyBook's own selection
# (A2) untouched so only the "active sheet" flag moves.
$wsInfo.Activate()
$wsInfo.Range("A4").Select()
